$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: Panasonic ERJ-1GNF1002C (10 kOhms resistor) - add R9 to the reference
# designators and bump the quantity from 2 to 3.
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = "R5, R9 R12"

# Row 25: Panasonic ERT-JZEG103FA (NTC Thermistor) - remove the TH2 reference
# designator and drop the quantity from 2 to 1.
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = "TH1"

# Update the worksheet's active selection / scroll position.
$ws.Range("C19").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1

# Update the workbook window size/position.
$win.Left = 3000
$win.Top = 690
$win.Width = 21600
$win.Height = 12735
